# Apply updated dSF (column F) values to Sheet1.
# These values represent re-pulled source data (per commit message:
# "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 4
    3  = 6
    4  = 1
    5  = -2
    6  = -1
    7  = 2
    8  = 6
    10 = 2
    11 = 10
    13 = -2
    14 = 1
    15 = 4
    16 = 1
    17 = 6
    18 = -2
    19 = 6
    20 = 1
    21 = 6
    22 = -1
    23 = 1
    24 = 2
    26 = -2
    27 = 6
    28 = -2
    29 = -3
    31 = 4
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
